$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Insert a new row above row 21 (shifts rows 21..111 down to 22..112,
# copying formatting from the row above per Excel's default behavior)
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 following the same pattern as the
# other data rows in this table
$ws.Cells.Item(21, 1).Value = "CREATE/MODIFY"
$ws.Cells.Item(21, 2).Value = "LIB_EWS_BE"
$ws.Cells.Item(21, 3).Value = "INDICATOR_215"
$ws.Cells.Item(21, 5).Value = "String"

# Update the selection to match the authored state
$ws.Range("C21").Select()
